# Update the LR-pairs sheet (Gm13306-Ackr2) with refreshed TPM-derived values.
# Adds a new "ECs" sending-cluster row, renames "Resolving-Mac" to "Neutrophils",
# and refreshes all detection/expression/specificity figures for rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster = ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gm13306"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.169054
$ws.Range("H2").Value = 0.338108
$ws.Range("I2").Value = 0.2998784017861604
$ws.Range("J2").Value = 0.2245827615387374
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2786473333333334
$ws.Range("N2").Value = 0.8359420000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.04710644628933334
$ws.Range("R2").Value = 0.282638677736
$ws.Range("S2").Value = 0.2998784017861604
$ws.Range("T2").Value = 0.2245827615387374

# Row 3: Sending cluster = FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gm13306"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3036856666666667
$ws.Range("H3").Value = 0.911057
$ws.Range("I3").Value = 0.5386963477076239
$ws.Range("J3").Value = 0.605154852825717
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2786473333333334
$ws.Range("N3").Value = 0.8359420000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.08462120118822224
$ws.Range("R3").Value = 0.761590810694
$ws.Range("S3").Value = 0.5386963477076239
$ws.Range("T3").Value = 0.605154852825717

# Row 4: Sending cluster = MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gm13306"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.0166775
$ws.Range("H4").Value = 0.033355
$ws.Range("I4").Value = 0.02958357711612082
$ws.Range("J4").Value = 0.02215551838798428
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2786473333333334
$ws.Range("N4").Value = 0.8359420000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.004647140901666668
$ws.Range("R4").Value = 0.02788284541000001
$ws.Range("S4").Value = 0.02958357711612082
$ws.Range("T4").Value = 0.02215551838798428

# Row 5: Sending cluster = Neutrophils
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Gm13306"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07432466666666666
$ws.Range("H5").Value = 0.222974
$ws.Range("I5").Value = 0.1318416733900949
$ws.Range("J5").Value = 0.1481068672475613
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2786473333333334
$ws.Range("N5").Value = 0.8359420000000001
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.02071037016755556
$ws.Range("R5").Value = 0.186393331508
$ws.Range("S5").Value = 0.1318416733900949
$ws.Range("T5").Value = 0.1481068672475613
